$d = $word.ActiveDocument

# Locate the existing "cool text and more fun." paragraph.
$r = $d.Content
$r.Find.Execute("cool text and more fun.", $false, $false, $false, $false,
                 $false, $true, 1, $false, "", 0) | Out-Null

# Split off a brand-new paragraph right after it containing the new sentence.
$r.InsertAfter("`r" + "just to make sure it works a 2nd time.")

# $r now spans both the original paragraph and the freshly-created one;
# the 2nd paragraph in that range is the new one - give it the BodyText style.
$newPara = $r.Paragraphs(2)
$newPara.Style = "BodyText"
